$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range("D26").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E26").Value = "['SoftwareFault']"

# Row 27
$ws.Range("D27").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['SoftwareFault']"

# Row 39
$ws.Range("D39").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E39").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

# Row 58
$ws.Range("D58").Value = "[0, 0, 0, 1, 0, 0, 1]"
$ws.Range("E58").Value = "['ParamViolation', 'SoftwareFault']"

# Row 73
$ws.Range("D73").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal']"

# Row 82
$ws.Range("D82").Value = "[1, 1, 1, 0, 0, 0, 0]"
$ws.Range("E82").Value = "['Normal', 'SurroundingEnvironment', 'HardwareFault']"
